$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, shifting existing rows 101-149 down to 102-150
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record's data
$ws.Cells.Item(101, 1).Value = 6
$ws.Cells.Item(101, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(101, 3).Value = "Metropolitana"
$ws.Cells.Item(101, 4).Value = 44489
$ws.Cells.Item(101, 5).Value = 13
$ws.Cells.Item(101, 6).Value = 100112022
$ws.Cells.Item(101, 7).Value = "Arveja Verde"
$ws.Cells.Item(101, 8).Value = "Perfection"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 300
$ws.Cells.Item(101, 11).Value = 20000
$ws.Cells.Item(101, 12).Value = 22000
$ws.Cells.Item(101, 13).Value = 21200
$ws.Cells.Item(101, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(101, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(101, 16).Value = 848
$ws.Cells.Item(101, 17).Value = 25
$ws.Cells.Item(101, 18).Value = "Hortaliza"
